# Auto-generated Excel COM-interop script implementing the pilot1_subjects.xlsx edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = 'Subjects'

# 2. Fix up defined names so they reference the new sheet name explicitly
foreach ($n in $wb.Names) {
    if ($n.Name -like '*_FilterDatabase*') {
        $n.RefersTo = '=Subjects!#REF!'
    }
    if ($n.Name -like '*Print_Titles*') {
        $n.RefersTo = '=Subjects!$1:$1'
    }
    if ($n.Name -eq 'valHighlight') {
        $n.RefersTo = '=IFERROR(IF(Subjects!#REF!="Sì", TRUE, FALSE),FALSE)'
    }
}

# 3. Update commentary cells: text + bold 'label' prefixes via Characters()
$ws.Range('I2').Value = 'TASK: not clear that the shuffles were re-starting from zero at every couple of graphs (he thought that they were 20 for each couple of graphs) [Safari version 16.4]'
$ws.Range('I2').Characters(1, 4).Font.Bold = $true

$ws.Range('I3').Value = 'TUTORIAL (I part): bad quality images during tutorial; TASK: thought that using more shuffles resulted in more difficult following trials'
$ws.Range('I3').Characters(1, 17).Font.Bold = $true
$ws.Range('I3').Characters(56, 5).Font.Bold = $true

$ws.Range('I4').Value = 'TUTORIAL (I part): when starting tutorial, she was expecting to see the red version first (it is the first time the black one is shown first)  + BROWSER COMPATIBILITY: problems visualizing the text (FillText) with Safari (version 14.1.1 ) + it''s written "right or left" arrow, but instead it should be "left or right" + she did not feel tired at all, could have done other 2/3 blocks'
$ws.Range('I4').Characters(1, 17).Font.Bold = $true

$ws.Range('I5').Value = 'INSTRUCTIONS: 1. indicate which is the last page when it is possible to you can go back, before moving on with tutorial; 2. Avoid repeating same senteces over and over, since it is possible to go back; 3. Convey the idea that the task will become very difficult, with some difficult examples in the tutorial.  TASK: communicate score of last trial of block. '
$ws.Range('I5').Characters(1, 12).Font.Bold = $true
$ws.Range('I5').Characters(311, 6).Font.Bold = $true

$ws.Range('I6').Value = 'INTRO PAGE: change title (New Version of...) + you will receive INSTRUCTIONS…  + in case you HAVE QUESTIONS regarding + AND THAT of the other volunteers. BROWSER compatibility: (Firefox -> informed consent page not fitting, is cut. Also score is cut during the task). TUTORIAL (II part): low quality of images. TASK: she was expecting to see the score after pressing spacebar, like in the tutorial (she suggested to make it more consistent in the two cases. Maybe by inverting the score increase - solution images). '
$ws.Range('I6').Characters(1, 10).Font.Bold = $true
$ws.Range('I6').Characters(268, 19).Font.Bold = $true
$ws.Range('I6').Characters(312, 4).Font.Bold = $true

$ws.Range('I7').Value = 'TASK: it was not clear what shuffles meant, he was thinking that one space bar press randomized the graph, and another one was taking the graph back to the previous visualization. Was very fast, but could not handle one trial more. Maybe he misunderstood / did not read the instructions carefully enough.'
$ws.Range('I7').Characters(1, 6).Font.Bold = $true

$ws.Range('I8').Value = 'TASK: she tried to be as fast as accurate and possible with no help, so she did not think of using the shuffles when she was unsure about the answer. Maybe one could underline more that there is no penalization in using the shuffles and that there is no time limit, and that shuffling helps in giving an answer. She could have handled 2 more blocks. NOTE: the internet connection was interrupted, so she had to restart the experiment. This is the reason for the fast execution of tutorial.'
$ws.Range('I8').Characters(1, 5).Font.Bold = $true

$ws.Range('I9').Value = 'INSTRUCTIONS: she had a feeling that at a point she got the point, the instructions could have been a bit shorter, especially in the part about the shuffles. Also at the beginning, maybe it was not necessary to move from the squares to triangles (going directly to explain the shuffles)'
$ws.Range('I9').Characters(1, 13).Font.Bold = $true

$ws.Range('I10').Value = 'INSTRUCTIONS: on Safari, "INSTRUCTIONS" is not centered in the screen. Possibility: using the same color for highlighted rows and a different one for columns? TASK: She was expecting to see red tiles appear after giving an answer'
$ws.Range('I10').Characters(1, 13).Font.Bold = $true
$ws.Range('I10').Characters(160, 6).Font.Bold = $true

$ws.Range('I11').Value = 'INTRO PAGE: privacy (servers located in EU = the ones of SISSA?) + also in case results WILL/WOULD BE PUBLISHED. Data will be handled according to EU REGULATION. Space after "GDPR". Dot after following sentence. INSTRUCTIONS: 1. using the same color for highlighted rows and a different one for columns? 2. Square-triangle transition: let''s cover A PART (not a half) of the square. 3. Using "..." between shaded and fully covered part. 4. A dot is missing in this part. 4. Switching two rows and columns TRANSFORMS the original triangle. 5. Vertical triangle (pointing right for the first time) -> let''s consider a FLIPPED VERSION of the inital triangle. 6. Ricontrolla di aver usato "regular checkerboard" o "chessboard" quando sono scacchiere regolari. 7. (shuffle in non-regular checkerboard) ALSO IN THIS CASE, a shuffle is the... 8. It is also possible to shuffle TRIANGLES WITH MORE TILES (not bigger). 9. ALSO IN THE CASE OF MORE TILES; a shuffle is the... 10. Lower the number of examples of shuffles with red clique (to make the tutorial faster)? 11. THE TRIANGLES OF THE EXPERIMENT WILL LOOK LIKE THIS (not "this will be the size of the bigger triangle...") 12. "One on the left... and one on the right" -> non far sparire il triangolo di sx quando si mostra quello di destra (risparmio anche di un''immagine. Se possibile, cerca anche di far apparire i cambiamenti sui triangoli allineati) TUTORIAL (pt. 2): "shuffles are limited" is repeated twice. 2. Organizza come: SHUFFLE 1 - SHUFFLE 2 - SHUFFLE 3 - FEEDBACK + SOLUZIONE (insieme, in modo che la visualizzazione successiva sia il trial successivo). 3. Sottolineare il fatto che nel task ci saranno più trials (dire proprio il numero -> "you will have ... shuffles for each couple of graphs"). FINAL CONSENT: "consent" or "consense"? Centrare l''elenco, ma non mettere i punti su linee diverse. TASK: 1. suggest to wear glasses at the beginning? 2. Feedback between blocks -> elenco puntato non centrato. 3. Decrease number of shuffles 4. Inter-trial interval più lungo per evitare sensory memory/after effect. '
$ws.Range('I11').Characters(1, 10).Font.Bold = $true
$ws.Range('I11').Characters(212, 14).Font.Bold = $true
$ws.Range('I11').Characters(1401, 17).Font.Bold = $true
$ws.Range('I11').Characters(1859, 5).Font.Bold = $true

# 4. Add new row 14 with the KEY POINT comment (copy I13's style, then set text)
$ws.Range('I13').Copy($ws.Range('I14'))
$ws.Range('I14').Value = 'KEY POINT: instruction information is not really used by the participants'
$ws.Rows(14).RowHeight = 24

# 5. Row 11 grew taller to fit the new wording
$ws.Rows(11).RowHeight = 201.6

# 6. Sheet view state: scroll position, zoom, active selection
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range('I3').Select()
